# Weekly update: insert a new price record for Agrícola del Norte S.A. de
# Arica - Durazno (Florida King, Región de O'Higgins) as row 52, pushing
# the existing rows 52-58 down to 53-59.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 52..58 down one position to make room for the new record.
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new weekly record.
$ws.Range("A52").Value = 1
$ws.Range("B52").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C52").Value = "Arica y Parinacota"
$ws.Range("D52").Value = 44918
$ws.Range("E52").Value = 15
$ws.Range("F52").Value = "Fruta"
$ws.Range("G52").Value = 100103
$ws.Range("H52").Value = "Frutos de hueso (carozo)"
$ws.Range("I52").Value = 100103004
$ws.Range("J52").Value = "Durazno"
$ws.Range("K52").Value = "Florida King"
$ws.Range("L52").Value = "Segunda"
$ws.Range("M52").Value = 450
$ws.Range("N52").Value = 13000
$ws.Range("O52").Value = 15000
$ws.Range("P52").Value = 14111
$ws.Range("Q52").Value = "`$/caja 18 kilos granel"
$ws.Range("R52").Value = "Región de O'Higgins"
$ws.Range("S52").Value = 784
$ws.Range("T52").Value = 18
